# Apply updated cryptocurrency price/volume data (and two row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "34.233.52"
Set-TextValue "E2" "  +0.44%  "

Set-TextValue "D3" "1.788.14"
Set-TextValue "E3" "  -0.02%  "

Set-TextValue "E4" "  +0.08%  "

Set-TextValue "D5" "225.90"
Set-TextValue "E5" "  -0.30%  "

Set-TextValue "E6" "  +0.49%  "

Set-TextValue "E7" "  -0.03%  "

Set-TextValue "D8" "32.28"
Set-TextValue "E8" "  +0.23%  "

Set-TextValue "E9" "  +0.23%  "

Set-TextValue "D10" "0.0690"
Set-TextValue "E10" "  +0.14%  "

Set-TextValue "E11" "  +0.81%  "

Set-TextValue "D12" "2.046.61"
Set-TextValue "E12" "  -0.01%  "

Set-TextValue "B13" "WrappedEther"
Set-TextValue "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.811.79"
Set-TextValue "E13" "  +0.75%  "

Set-TextValue "B14" "Chainlink"
Set-TextValue "C14" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D14" "11.14"
Set-TextValue "E14" "  -3.14%  "

Set-TextValue "D15" "0.626"
Set-TextValue "E15" "  +0.73%  "

Set-TextValue "D16" "34.223.70"
Set-TextValue "E16" "  +0.41%  "

Set-TextValue "E17" "  +0.37%  "

Set-TextValue "D18" "67.98"
Set-TextValue "E18" "  +0.02%  "

Set-TextValue "D19" "0.0₃0807"
Set-TextValue "E19" "  +3.52%  "

Set-TextValue "D20" "246.09"
Set-TextValue "E20" "  +0.92%  "

Set-TextValue "D21" "10.99"
Set-TextValue "E21" "  +0.75%  "

Set-TextValue "D22" "0.999"
Set-TextValue "E22" "  -0.05%  "

Set-TextValue "E23" "  +1.53%  "

Set-TextValue "E24" "  +0.77%  "

Set-TextValue "D25" "161.70"
Set-TextValue "E25" "  -0.33%  "

Set-TextValue "E26" "  -0.25%  "

Set-TextValue "D27" "16.34"
Set-TextValue "E27" "  +0.27%  "

Set-TextValue "E28" "  +0.81%  "

Set-TextValue "E29" "  +0.21%  "

Set-TextValue "B30" "Hedera"
Set-TextValue "C30" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D30" "0.0521"
Set-TextValue "E30" "  +0.08%  "

Set-TextValue "B31" "PancakeSwap"
Set-TextValue "C31" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "1.23"
Set-TextValue "E31" "  -0.30%  "

Set-TextValue "D32" "3.77"
Set-TextValue "E32" "  +3.03%  "

Set-TextValue "E33" "  +3.99%  "

Set-TextValue "E34" "  -1.20%  "

Set-TextValue "D35" "1.443.93"
Set-TextValue "E35" "  +2.09%  "

Set-TextValue "D36" "2.57"
Set-TextValue "E36" "  +9.27%  "

Set-TextValue "D37" "0.664"
Set-TextValue "E37" "  +2.90%  "

Set-TextValue "E38" "  +1.27%  "

Set-TextValue "E39" "  -0.19%  "

Set-TextValue "D40" "82.14"
Set-TextValue "E40" "  +2.10%  "

Set-TextValue "D41" "2.38"
Set-TextValue "E41" "  +1.06%  "

Set-TextValue "D42" "14.03"
Set-TextValue "E42" "  +4.21%  "

Set-TextValue "E43" "  +0.27%  "

Set-TextValue "D44" "2.71"
Set-TextValue "E44" "  +1.13%  "

Set-TextValue "E45" "  +2.44%  "

Set-TextValue "E46" "  +0.28%  "

Set-TextValue "E47" "  +0.65%  "

Set-TextValue "D48" "1.946.89"
Set-TextValue "E48" "  -0.10%  "

Set-TextValue "E49" "  -1.60%  "

Set-TextValue "E50" "  -0.01%  "

Set-TextValue "D51" "0.0₆0126"
Set-TextValue "E51" "  -8.09%  "
